$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (old) values of B2:F10 before overwriting anything,
# since rows 3-11 will receive the old values of rows 2-10 (shift down by one row).
$oldValues = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @()
    for ($c = 2; $c -le 6; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value()
    }
    $oldValues[$r] = $rowVals
}

# Shift old rows 2-10 down into rows 3-11
for ($r = 10; $r -ge 2; $r--) {
    $vals = $oldValues[$r]
    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r + 1, $c).Value = $vals[$c - 2]
    }
}

# Write the new values into row 2
$ws.Range("B2").Value = 0.1560865643779764
$ws.Range("C2").Value = 0.3515450347245845
$ws.Range("D2").Value = 0.2150214299408537
$ws.Range("E2").Value = 0.4637040326985023
$ws.Range("F2").Value = 0.446457955381491
